$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - match the bold/bordered header style used by
# the existing header row (A1:L1) by copying the format from L1.
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

# Fill the new data columns for every existing data row (2-13)
$lastRow = 13
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20120398
    $ws.Cells.Item($r, 15).Value = 1
}
